# fix the merge error
# Player.xlsx -> "Property" sheet: the View (column F) flag was dropped for
# several rows during a bad merge, and the Private/Save/View values for the
# GameID/GateID rows got shifted by one column. GuildID's Public flag was
# also incorrectly left as TRUE. Restore the intended values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

# Rows 68-75 (X, Y, Z, TargetX, TargetY, TargetZ, PathStep, LoadPropertyFinish):
# the View column (F) was blank and should be TRUE.
for ($r = 68; $r -le 75; $r++) {
    $ws.Cells.Item($r, 6).Value = $true
}

# Row 76 (GameID): Private/Save/View values were shifted one column to the
# left during the merge (D/E held TRUE, TRUE, F held FALSE). Correct values
# are Private=FALSE, Save=FALSE, View=TRUE.
$ws.Cells.Item(76, 4).Value = $false
$ws.Cells.Item(76, 5).Value = $false
$ws.Cells.Item(76, 6).Value = $true

# Row 77 (GateID): same fix as row 76.
$ws.Cells.Item(77, 4).Value = $false
$ws.Cells.Item(77, 5).Value = $false
$ws.Cells.Item(77, 6).Value = $true

# Row 78 (GuildID): Public flag incorrectly left TRUE, should be FALSE.
$ws.Cells.Item(78, 3).Value = $false

# Restore the view/selection state that was active when the fix was made.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C78").Select() | Out-Null

$wb.Save()
